$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above the current row 11, shifting the existing
# rows 11-23 down to rows 13-25 (keeping their data intact).
$ws.Rows("11:12").Insert()

# --- New row 11 ---
$ws.Range("A11").Value2 = 2
$ws.Range("B11").Value2 = "Comercializadora del Agro de Limarí"
$ws.Range("C11").Value2 = "Coquimbo"
$ws.Range("D11").Value2 = 44580
$ws.Range("E11").Value2 = 4
$ws.Range("F11").Value2 = "Fruta"
$ws.Range("G11").Value2 = 100103
$ws.Range("H11").Value2 = "Frutos de hueso (carozo)"
$ws.Range("I11").Value2 = 100103003
$ws.Range("J11").Value2 = "Damasco"
$ws.Range("K11").Value2 = "Modesto"
$ws.Range("L11").Value2 = "Especial"
$ws.Range("M11").Value2 = 300
$ws.Range("N11").Value2 = 22500
$ws.Range("O11").Value2 = 23000
$ws.Range("P11").Value2 = 22750
$ws.Range("Q11").Value2 = "$/caja 18 kilos"
$ws.Range("R11").Value2 = "Región Metropolitana"
$ws.Range("S11").Value2 = 1264
$ws.Range("T11").Value2 = 18

# --- New row 12 ---
$ws.Range("A12").Value2 = 2
$ws.Range("B12").Value2 = "Comercializadora del Agro de Limarí"
$ws.Range("C12").Value2 = "Coquimbo"
$ws.Range("D12").Value2 = 44580
$ws.Range("E12").Value2 = 4
$ws.Range("F12").Value2 = "Fruta"
$ws.Range("G12").Value2 = 100103
$ws.Range("H12").Value2 = "Frutos de hueso (carozo)"
$ws.Range("I12").Value2 = 100103003
$ws.Range("J12").Value2 = "Damasco"
$ws.Range("K12").Value2 = "Modesto"
$ws.Range("L12").Value2 = "Primera"
$ws.Range("M12").Value2 = 400
$ws.Range("N12").Value2 = 19500
$ws.Range("O12").Value2 = 20000
$ws.Range("P12").Value2 = 19750
$ws.Range("Q12").Value2 = "$/caja 18 kilos"
$ws.Range("R12").Value2 = "Región Metropolitana"
$ws.Range("S12").Value2 = 1097
$ws.Range("T12").Value2 = 18
